$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 15 (@@ -1346,22 +1346,22 @@)
$ws.Range("H15").Value = 651.2033699999999
$ws.Range("I15").Value = 651.2033699999999
$ws.Range("K15").Value = 1953.61011
$ws.Range("M15").Value = -1784.61011
# row 40 (@@ -2580,25 +2580,25 @@)
$ws.Range("H40").Value = 1303.75
$ws.Range("J40").Value = 1646
$ws.Range("L40").Value = 1646
$ws.Range("N40").Value = -1996
# row 62 (@@ -3673,25 +3673,25 @@)
$ws.Range("H62").Value = 4721.222
$ws.Range("J62").Value = 6331.8335
$ws.Range("L62").Value = 6331.8335
$ws.Range("N62").Value = -7579.8335
# row 65 (@@ -3823,25 +3823,25 @@)
$ws.Range("H65").Value = 4721.222
$ws.Range("J65").Value = 6331.8335
$ws.Range("L65").Value = 31659.1675
$ws.Range("N65").Value = -37899.1675
# row 80 (@@ -4555,25 +4555,25 @@)
$ws.Range("H80").Value = 2317.5
$ws.Range("I80").Value = 3967.3333
$ws.Range("J80").Value = 667.6667
$ws.Range("K80").Value = 11901.9999
$ws.Range("L80").Value = 2003.0001
$ws.Range("M80").Value = -10903.9999
$ws.Range("N80").Value = -3999.0001
# row 83 (@@ -4705,25 +4705,25 @@)
$ws.Range("H83").Value = 2317.5
$ws.Range("I83").Value = 3967.3333
$ws.Range("J83").Value = 667.6667
$ws.Range("K83").Value = 35705.9997
$ws.Range("L83").Value = 6009.0003
$ws.Range("M83").Value = -30713.9997
$ws.Range("N83").Value = -15993.0003
# row 107 (@@ -5917,25 +5917,25 @@)
$ws.Range("H107").Value = 789
$ws.Range("I107").Value = 570
$ws.Range("J107").Value = 1095.6
$ws.Range("K107").Value = 570
$ws.Range("L107").Value = 1095.6
$ws.Range("M107").Value = 1350
$ws.Range("N107").Value = -4935.6
# row 132 (@@ -7145,25 +7145,25 @@)
$ws.Range("H132").Value = 731.6667
$ws.Range("I132").Value = 738.4
$ws.Range("J132").Value = 664.3333
$ws.Range("K132").Value = 2215.2
$ws.Range("L132").Value = 1992.9999
$ws.Range("M132").Value = 314.8000000000002
$ws.Range("N132").Value = -7052.9999
# row 137 (@@ -7393,22 +7393,22 @@)
$ws.Range("H137").Value = 1826.4333
$ws.Range("I137").Value = 1150.8
$ws.Range("K137").Value = 3452.4
$ws.Range("M137").Value = -902.3999999999996
# row 138 (@@ -7445,25 +7445,22 @@)
$ws.Range("H138").Value = 1532.9
$ws.Range("I138").Value = 1532.9
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 4598.700000000001
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 541.2999999999993
$ws.Range("N138").ClearContents()

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32 (@@ -9187,22 +9184,22 @@)
$ws.Range("H32").Value = 4519.436
$ws.Range("I32").Value = 3081.9333
$ws.Range("K32").Value = 3081.9333
$ws.Range("M32").Value = -2794.9333
# row 45 (@@ -9818,22 +9815,22 @@)
$ws.Range("H45").Value = 3464198.2
$ws.Range("I45").Value = 11252559
$ws.Range("K45").Value = 11252559
$ws.Range("M45").Value = -11252182
# row 61 (@@ -10584,25 +10581,25 @@)
$ws.Range("H61").Value = 3200.353
$ws.Range("I61").Value = 2230.182
$ws.Range("J61").Value = 4979
$ws.Range("K61").Value = 2230.182
$ws.Range("L61").Value = 4979
$ws.Range("M61").Value = -2018.182
$ws.Range("N61").Value = -5403
# row 74 (@@ -11206,22 +11203,22 @@)
$ws.Range("H74").Value = 1209.2916
$ws.Range("I74").Value = 572.7222
$ws.Range("K74").Value = 572.7222
$ws.Range("M74").Value = 301.2778
# row 77 (@@ -11350,22 +11347,22 @@)
$ws.Range("H77").Value = 1209.2916
$ws.Range("I77").Value = 572.7222
$ws.Range("K77").Value = 2863.611
$ws.Range("M77").Value = 1504.389
# row 136 (@@ -14181,25 +14178,25 @@)
$ws.Range("H136").Value = 3200.353
$ws.Range("I136").Value = 2230.182
$ws.Range("J136").Value = 4979
$ws.Range("K136").Value = 6690.545999999999
$ws.Range("L136").Value = 14937
$ws.Range("M136").Value = -4140.545999999999
$ws.Range("N136").Value = -20037

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 82 (@@ -18372,25 +18369,25 @@)
$ws.Range("H82").Value = 42559.8
$ws.Range("J82").Value = 44933.332
$ws.Range("L82").Value = 44933.332
$ws.Range("N82").Value = -45699.332
# row 85 (@@ -18525,25 +18522,25 @@)
$ws.Range("H85").Value = 42559.8
$ws.Range("J85").Value = 44933.332
$ws.Range("L85").Value = 44933.332
$ws.Range("N85").Value = -47585.332
# row 105 (@@ -19496,25 +19493,25 @@)
$ws.Range("H105").Value = 2519.9
$ws.Range("I105").Value = 2538.7778
$ws.Range("J105").Value = 2350
$ws.Range("K105").Value = 2538.7778
$ws.Range("L105").Value = 2350
$ws.Range("M105").Value = -791.7777999999998
$ws.Range("N105").Value = -5844

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31 (@@ -22764,25 +22761,25 @@)
$ws.Range("H31").Value = 1440.6154
$ws.Range("I31").Value = 978.4
$ws.Range("J31").Value = 1550.6666
$ws.Range("K31").Value = 978.4
$ws.Range("L31").Value = 1550.6666
$ws.Range("M31").Value = -683.4
$ws.Range("N31").Value = -2140.6666
# row 34 (@@ -22917,25 +22914,25 @@)
$ws.Range("H34").Value = 1440.6154
$ws.Range("I34").Value = 978.4
$ws.Range("J34").Value = 1550.6666
$ws.Range("K34").Value = 978.4
$ws.Range("L34").Value = 1550.6666
$ws.Range("M34").Value = -776.4
$ws.Range("N34").Value = -1954.6666
# row 105 (@@ -26366,7 +26363,7 @@)
$ws.Range("H105").Value = 620
# row 110 (@@ -26611,19 +26608,22 @@)
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180
# row 134 (@@ -27751,22 +27751,22 @@)
$ws.Range("H134").Value = 1304.9667
$ws.Range("I134").Value = 1029.6296
$ws.Range("K134").Value = 3088.8888
$ws.Range("M134").Value = -553.8887999999997

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 68 (@@ -31567,25 +31567,25 @@)
$ws.Range("H68").Value = 1599.8864
$ws.Range("I68").Value = 770.44446
$ws.Range("J68").Value = 1813.1714
$ws.Range("K68").Value = 2311.33338
$ws.Range("L68").Value = 5439.5142
$ws.Range("M68").Value = -1500.33338
$ws.Range("N68").Value = -7061.5142
# row 71 (@@ -31717,25 +31717,25 @@)
$ws.Range("H71").Value = 1599.8864
$ws.Range("I71").Value = 770.44446
$ws.Range("J71").Value = 1813.1714
$ws.Range("K71").Value = 6934.00014
$ws.Range("L71").Value = 16318.5426
$ws.Range("M71").Value = -2878.00014
$ws.Range("N71").Value = -24430.5426
# row 92 (@@ -32773,25 +32773,25 @@)
$ws.Range("H92").Value = 614.5714
$ws.Range("J92").Value = 667
$ws.Range("L92").Value = 2001
$ws.Range("N92").Value = -4497
# row 115 (@@ -33933,25 +33933,22 @@)
$ws.Range("H115").Value = 1599.6666
$ws.Range("I115").Value = 1599.6666
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 4798.9998
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -3623.9998
$ws.Range("N115").ClearContents()
# row 119 (@@ -34138,22 +34135,25 @@)
$ws.Range("H119").Value = 62503750
$ws.Range("I119").Value = 250000000
$ws.Range("J119").Value = 5000
$ws.Range("K119").Value = 750000000
$ws.Range("L119").Value = 15000
$ws.Range("M119").Value = -749995162
$ws.Range("N119").Value = -24676
# row 122 (@@ -34288,22 +34288,22 @@)
$ws.Range("H122").Value = 1562.6666
$ws.Range("I122").Value = 1344
$ws.Range("K122").Value = 12096
$ws.Range("M122").Value = -9646

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 132 (@@ -41696,25 +41696,25 @@)
$ws.Range("H132").Value = 1833752.8
$ws.Range("I132").Value = 2566153.2
$ws.Range("J132").Value = 2751.5
$ws.Range("K132").Value = 7698459.600000001
$ws.Range("L132").Value = 8254.5
$ws.Range("M132").Value = -7695929.600000001
$ws.Range("N132").Value = -13314.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7 (@@ -42519,25 +42519,25 @@)
$ws.Range("H7").Value = 4387.2144
$ws.Range("J7").Value = 3955.1538
$ws.Range("L7").Value = 3955.1538
$ws.Range("N7").Value = -4179.1538
# row 38 (@@ -44029,19 +44029,22 @@)
$ws.Range("H38").Value = 10000
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10820
# row 61 (@@ -45129,25 +45132,25 @@)
$ws.Range("H61").Value = 3285.5293
$ws.Range("J61").Value = 3337.8333
$ws.Range("L61").Value = 3337.8333
$ws.Range("N61").Value = -3741.8333
# row 93 (@@ -46673,22 +46676,22 @@)
$ws.Range("H93").Value = 1499.25
$ws.Range("I93").Value = 998.5
$ws.Range("K93").Value = 998.5
$ws.Range("M93").Value = 249.5
# row 113 (@@ -47635,25 +47638,25 @@)
$ws.Range("H113").Value = 3285.5293
$ws.Range("J113").Value = 3337.8333
$ws.Range("L113").Value = 3337.8333
$ws.Range("N113").Value = -7677.8333
# row 122 (@@ -48055,22 +48058,22 @@)
$ws.Range("H122").Value = 8305.786
$ws.Range("I122").Value = 6785.875
$ws.Range("K122").Value = 20357.625
$ws.Range("M122").Value = -17907.625
# row 126 (@@ -48248,25 +48251,25 @@)
$ws.Range("H126").Value = 4387.2144
$ws.Range("J126").Value = 3955.1538
$ws.Range("L126").Value = 11865.4614
$ws.Range("N126").Value = -16805.4614

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 126 (@@ -55097,25 +55100,25 @@)
$ws.Range("H126").Value = 15344.363
$ws.Range("I126").Value = 21480.5
$ws.Range("J126").Value = 7981
$ws.Range("K126").Value = 64441.5
$ws.Range("L126").Value = 23943
$ws.Range("M126").Value = -61971.5
$ws.Range("N126").Value = -28883
# row 132 (@@ -55385,25 +55388,25 @@)
$ws.Range("H132").Value = 1412.8837
$ws.Range("I132").Value = 1192.6061
$ws.Range("J132").Value = 2139.8
$ws.Range("K132").Value = 3577.8183
$ws.Range("L132").Value = 6419.400000000001
$ws.Range("M132").Value = -1047.8183
$ws.Range("N132").Value = -11479.4
